$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.356.62'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.339.14'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -4.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -5.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.16'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.58%  '
$ws.Range("E7").Value = '  -2.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.329.90'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.15%  '
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.162'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.67'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.33%  '
$ws.Range("E13").Value = '  -2.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.04'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.874.93'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.29'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("E17").Value = '  -3.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.337.72'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -4.29%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '64.263.07'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.73'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.978'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '432.89'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.79%  '
$ws.Range("E23").Value = '  +16.22%  '
$ws.Range("E24").Value = '  -6.10%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.29'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.43'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.73'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.82'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.72'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.71'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.66'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("E32").Value = '  -2.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '579.59'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.57%  '
$ws.Range("E34").Value = '  -3.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.26'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.44%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  -8.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.48'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.58'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0751'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.366'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.105.06'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.80'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.23'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.72%  '
$ws.Range("E46").Value = '  -2.82%  '
$ws.Range("E47").Value = '  -3.73%  '
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.59'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.28'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '134.76'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.73%  '
